$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.325.22"
$ws.Range("E2").Value = "  -0.99%  "
$ws.Range("D3").Value = "3.235.51"
$ws.Range("E3").Value = "  +2.75%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "595.40"
$ws.Range("E5").Value = "  -1.10%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.20"
$ws.Range("E6").Value = "  -1.20%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "3.231.79"
$ws.Range("E8").Value = "  +2.85%  "
$ws.Range("E9").Value = "  -1.70%  "
$ws.Range("E10").Value = "  -1.20%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.38"
$ws.Range("E11").Value = "  -0.10%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.465"
$ws.Range("E12").Value = "  -0.50%  "
$ws.Range("E13").Value = "  -2.90%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.35"
$ws.Range("E14").Value = "  -1.88%  "
$ws.Range("D15").Value = "3.765.82"
$ws.Range("E15").Value = "  +2.73%  "
$ws.Range("E16").Value = "  +0.09%  "
$ws.Range("D17").Value = "3.230.92"
$ws.Range("E17").Value = "  +2.80%  "
$ws.Range("D18").Value = "63.334.88"
$ws.Range("E18").Value = "  -1.06%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.79"
$ws.Range("E19").Value = "  -1.09%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "473.41"
$ws.Range("E20").Value = "  -3.05%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.21"
$ws.Range("E21").Value = "  -3.27%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.730"
$ws.Range("E22").Value = "  +2.41%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.93"
$ws.Range("E23").Value = "  +2.31%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.74"
$ws.Range("E24").Value = "  -5.14%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.19"
$ws.Range("E25").Value = "  -0.81%  "
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.54"
$ws.Range("E27").Value = "  +7.62%  "
$ws.Range("E28").Value = "  -0.93%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.11"
$ws.Range("E29").Value = "  -1.14%  "
$ws.Range("E30").Value = "  +2.63%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "27.41"
$ws.Range("E31").Value = "  -1.15%  "
$ws.Range("E32").Value = "  -0.09%  "
$ws.Range("E33").Value = "  -4.40%  "
$ws.Range("E34").Value = "  -4.68%  "
$ws.Range("E35").Value = "  -1.54%  "
$ws.Range("E36").Value = "  -2.37%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "52.63"
$ws.Range("E37").Value = "  -0.09%  "
$ws.Range("E38").Value = "  -5.38%  "
$ws.Range("E39").Value = "  -1.23%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "422.34"
$ws.Range("E40").Value = "  -2.52%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.39"
$ws.Range("E41").Value = "  +0.22%  "
$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.76"
$ws.Range("E42").Value = "  -6.48%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "2.966.81"
$ws.Range("E43").Value = "  +1.15%  "
$ws.Range("E44").Value = "  -8.86%  "
$ws.Range("E45").Value = "  +2.67%  "
$ws.Range("E46").Value = "  -1.42%  "
$ws.Range("E47").Value = "  +0.06%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.35"
$ws.Range("E48").Value = "  -2.23%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "25.97"
$ws.Range("E49").Value = "  +0.41%  "
$ws.Range("E50").Value = "  -0.48%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "121.02"
$ws.Range("E51").Value = "  +0.29%  "
